# ------------------------------------------------------------------
# Applies the scraped "daily rebuild" update to 上海-漫展信息.xlsx
#   Sheet1 = 展览 (Exhibitions)
#   Sheet2 = 演出 (Performances)
#   Sheet3 = 本地生活 (Local life)
#   Sheet4 = 全部类型 (All types)
# Mostly "想去人数" (want-to-go counter, column F) bumps, one
# ticket-status cell (G31 on sheet1) flipping from "不可售" text to a
# real price, and a brand-new exhibition inserted at row 16 of sheet4
# which pushes the following 8 rows down by one.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===================== Sheet1 : 展览 =====================
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value  = 126
$ws1.Range("F3").Value  = 951
$ws1.Range("F4").Value  = 593
$ws1.Range("F5").Value  = 2832
$ws1.Range("F8").Value  = 592
$ws1.Range("F9").Value  = 392
$ws1.Range("F11").Value = 373
$ws1.Range("F12").Value = 440
$ws1.Range("F13").Value = 523
$ws1.Range("F14").Value = 2146
$ws1.Range("F15").Value = 1251
$ws1.Range("F16").Value = 730
$ws1.Range("F18").Value = 2657
$ws1.Range("F21").Value = 517
$ws1.Range("F22").Value = 526
$ws1.Range("F24").Value = 565
$ws1.Range("F25").Value = 565
$ws1.Range("F26").Value = 19
$ws1.Range("F28").Value = 548
$ws1.Range("F29").Value = 560
$ws1.Range("F30").Value = 223
$ws1.Range("F32").Value = 367
$ws1.Range("F34").Value = 215

# Row 31 ticket price becomes available (was inline text "不可售")
$ws1.Range("G31").Value = 89

# ===================== Sheet2 : 演出 =====================
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F7").Value  = 327
$ws2.Range("F21").Value = 272
$ws2.Range("F32").Value = 501
$ws2.Range("F33").Value = 501

# ===================== Sheet3 : 本地生活 =====================
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F4").Value = 1449
$ws3.Range("F6").Value = 199
$ws3.Range("F7").Value = 230

# ===================== Sheet4 : 全部类型 =====================
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value  = 1449
$ws4.Range("F4").Value  = 126
$ws4.Range("F6").Value  = 199
$ws4.Range("F8").Value  = 951
$ws4.Range("F9").Value  = 593
$ws4.Range("F10").Value = 2832
$ws4.Range("F11").Value = 2832
$ws4.Range("F14").Value = 592
$ws4.Range("F15").Value = 392

# A new event ("上海·第四届次元鹿角动漫游戏展") was scraped in at row 16,
# pushing the previously-16..24 rows down to 17..24 (with their "想去
# 人数" counters refreshed to match the values used elsewhere in this
# same export). Row 25 onward is unaffected (only a later counter bump).
$ws4.Range("B16").Value = "2024-04-05"
$ws4.Range("C16").Value = "上海·第四届次元鹿角动漫游戏展"
$ws4.Range("D16").Value = "长宁路1191号来福士西区(W)B1层01号、11号 星零界"
$ws4.Range("E16").Value = "2024.04.05 10:00-04.06 17:00"
$ws4.Range("F16").Value = 653
$ws4.Range("G16").Value = 68
$ws4.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=78228"
$ws4.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202403/0UZF1n651711425506347.jpeg"

$ws4.Range("B17").Value = "2024-04-08"
$ws4.Range("C17").Value = "上海·Walk Off The Earth 逃离地球2024巡演"
$ws4.Range("D17").Value = "汶水路116号 上海静安体育中心"
$ws4.Range("E17").Value = "2024.04.08 19:30-04.08 22:00"
$ws4.Range("F17").Value = 4
$ws4.Range("G17").Value = 580
$ws4.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=83418"
$ws4.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202403/l4HTKMoH1711520043516.jpeg"

$ws4.Range("B18").Value = "2024-04-12"
$ws4.Range("C18").Value = "上海·吉卜力工作室物语-沉浸式艺术展全球首站"
$ws4.Range("D18").Value = "龙台路10号2F 上海国际传媒港艺术中心"
$ws4.Range("E18").Value = "2024.04.12 10:00-05.12 20:00"
$ws4.Range("F18").Value = 440
$ws4.Range("G18").Value = 158
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=83036"
$ws4.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202403/aZoum5Hd1710472525792.jpeg"

$ws4.Range("B19").Value = "2024-04-12"
$ws4.Range("C19").Value = "上海·奇迹の闪耀 「UP!」巡回动漫演唱会"
$ws4.Range("D19").Value = "北京西路1700号 云峰剧院"
$ws4.Range("E19").Value = "2024.04.12 19:30-04.12 21:30"
$ws4.Range("F19").Value = 18
$ws4.Range("G19").Value = 126
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=82427"
$ws4.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202403/HvxHPz981709707512970.jpeg"

$ws4.Range("B20").Value = "2024-04-12"
$ws4.Range("C20").Value = "上海·铃木木乃美 2024 演唱会"
$ws4.Range("D20").Value = "宜昌路179号 万代南梦宫上海文化中心"
$ws4.Range("E20").Value = "2024.04.12 19:00-04.12 20:30"
$ws4.Range("F20").Value = 327
$ws4.Range("G20").Value = 380
$ws4.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=81906"
$ws4.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202402/rGjpHpAV1708328728461.jpeg"

$ws4.Range("B21").Value = "2024-04-13"
$ws4.Range("C21").Value = "上海·《四月是你的谎言》——“公生”与“薰”的钢琴小提琴唯美经典音乐集"
$ws4.Range("D21").Value = "丁香路425号 上海东方艺术中心"
$ws4.Range("E21").Value = "2024.04.13 19:30-04.13 21:30"
$ws4.Range("F21").Value = 345
$ws4.Range("G21").Value = 80
$ws4.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=78667"
$ws4.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202311/bTP7w6GD1700130122940.jpeg"

$ws4.Range("B22").Value = "2024-04-13"
$ws4.Range("C22").Value = "上海·【早鸟5折】红楼梦·梁祝·探清水河 《国潮》跨界音乐会"
$ws4.Range("D22").Value = "北海路251号近西藏南路 茉莉花剧场"
$ws4.Range("E22").Value = "2024.04.13 19:30-04.13 21:00"
$ws4.Range("F22").Value = 18
$ws4.Range("G22").Value = 90
$ws4.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=82790"
$ws4.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202403/CoBallQU1710311232127.jpeg"

$ws4.Range("B23").Value = "2024-04-13"
$ws4.Range("C23").Value = "上海·坏孩纸物语第38届动漫节之聂政篇"
$ws4.Range("D23").Value = "万航渡路889号 悦达889商业广场"
$ws4.Range("E23").Value = "2024.04.13 11:00-04.14 17:00"
$ws4.Range("F23").Value = 523
$ws4.Range("G23").Value = 63.9
$ws4.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=83266"
$ws4.Range("I23").Value = "//i1.hdslb.com/bfs/openplatform/202403/2TB6W9aP1711090258417.png"

$ws4.Range("B24").Value = "2024-04-13"
$ws4.Range("C24").Value = "上海·第三届奇卡波利国潮嘉年华-原X铁X崩同好交流"
$ws4.Range("D24").Value = "漕宝路3366号 七宝万科广场"
$ws4.Range("E24").Value = "2024.04.13 10:30-04.14 16:30"
$ws4.Range("F24").Value = 2146
$ws4.Range("G24").Value = 58.8
$ws4.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=82376"
$ws4.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202403/64i0bjSy1709692398951.jpeg"

# Counter-only bumps for the remainder of the sheet
$ws4.Range("F25").Value = 730
$ws4.Range("F28").Value = 2657
$ws4.Range("F30").Value = 520
$ws4.Range("F31").Value = 526
$ws4.Range("F33").Value = 230
$ws4.Range("F38").Value = 565
$ws4.Range("F39").Value = 565
$ws4.Range("F40").Value = 272
$ws4.Range("F41").Value = 548
$ws4.Range("F42").Value = 560
$ws4.Range("F45").Value = 367
$ws4.Range("F49").Value = 501
$ws4.Range("F50").Value = 501
